$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting the existing row 3 (and below) down
$ws.Rows("3:3").Insert()

# Fill in the newly inserted row 3
$ws.Range("A3").Value = 333
$ws.Range("B3").Value = 444
$ws.Range("C3").Value = "Fre"
$ws.Range("D3").Value = "Lip"

# Add a new row 5 with data
$ws.Range("A5").Value = 531
$ws.Range("B5").Value = 290
$ws.Range("C5").Value = "cq"
$ws.Range("D5").Value = "tk"
